$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 ---
$ws.Range("AI3").Value = 11
$ws.Range("AQ3").Value = 10

# --- Row 5 ---
$ws.Range("AB5").Value = $null
$ws.Range("AD5").Value = 4
$ws.Range("AF5").Value = 5
$ws.Range("AI5").Value = 15
$ws.Range("AK5").Value = 14
$ws.Range("AM5").Value = $null
$ws.Range("AM5").HorizontalAlignment = -4131
$ws.Range("AO5").Value = 13
$ws.Range("AQ5").Value = 12

# --- Row 6 ---
$ws.Range("Y6").Value = 13
$ws.Range("AA6").Value = 6
$ws.Range("AC6").Value = 7
$ws.Range("AE6").Value = 14
$ws.Range("AL6").Value = 17
$ws.Range("AN6").Value = 16

# --- Row 7 ---
$ws.Range("Z7").Value = 16
$ws.Range("AD7").Value = 17

# --- Row 10 ---
$ws.Range("AA10").Value = 8
$ws.Range("AC10").Value = 9
$ws.Range("AL10").Value = 19
$ws.Range("AN10").Value = 18

# --- Rows 15-16: clear Q/R/S columns ---
$ws.Range("Q15:S16").ClearContents()

# --- Rows 17-20: shift data up and renumber ---
$ws.Range("M17").Value = 3
$ws.Range("N17").Value = 6
$ws.Range("O17").Value = 7

$ws.Range("M18").Value = 3
$ws.Range("N18").Value = 7
$ws.Range("O18").Value = 4

$ws.Range("M19").Value = 8
$ws.Range("N19").Value = 9
$ws.Range("O19").Value = 7
$ws.Range("Q19:S19").ClearContents()

$ws.Range("M20").Value = 8
$ws.Range("N20").Value = 7
$ws.Range("O20").Value = 6
$ws.Range("Q20:S20").ClearContents()

# --- Rows 21-22: clear M value, keep style ---
$ws.Range("M21").Value = $null
$ws.Range("M22").Value = $null

# --- Sheet view: scroll back to top and select Q15 ---
$ws.Range("Q15").Select() | Out-Null
